$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "$ 27.333 CLP 19-10-20"
$ws.Range("A6").Value = "$ 27.333 CLP 19-10-20"
